# Rotate the artfynd (species find) records held in rows 5, 6 and 7:
#   new row5 <- old row7
#   new row6 <- old row5
#   new row7 <- old row6
# Only the columns that actually differ between the three rows are
# touched; shared columns (C, J, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AW, AX, ...) are identical across the three records and are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture "before" snapshots of the columns that vary between the
#     three rows (numeric columns) ---
$numCols = @("A","B","E","Q","R")
$row5 = @{}
$row6 = @{}
$row7 = @{}
foreach ($col in $numCols) {
    $row5[$col] = $ws.Range("$col`5").Value()
    $row6[$col] = $ws.Range("$col`6").Value()
    $row7[$col] = $ws.Range("$col`7").Value()
}

# --- text columns that vary between the three rows ---
$txtCols = @("D","F","G","H")
$row5t = @{}
$row6t = @{}
$row7t = @{}
foreach ($col in $txtCols) {
    $row5t[$col] = $ws.Range("$col`5").Value()
    $row6t[$col] = $ws.Range("$col`6").Value()
    $row7t[$col] = $ws.Range("$col`7").Value()
}

# "Antal" (I) is stored as text even though it looks numeric - capture its
# value and remember to write it back with a leading apostrophe so it does
# not get auto-coerced into a number.
$i5 = $ws.Range("I5").Value().ToString()
$i6 = $ws.Range("I6").Value().ToString()
$i7 = $ws.Range("I7").Value().ToString()

# "Publik kommentar" (AC) - present on rows 6 and 7 only (row 5 has none).
$ac5 = $ws.Range("AC5").Value()
$ac6 = $ws.Range("AC6").Value()
$ac7 = $ws.Range("AC7").Value()

# --- write the rotated values back: row5 <- row7, row6 <- row5, row7 <- row6 ---
foreach ($col in $numCols) {
    $ws.Range("$col`5").Value = $row7[$col]
    $ws.Range("$col`6").Value = $row5[$col]
    $ws.Range("$col`7").Value = $row6[$col]
}
foreach ($col in $txtCols) {
    $ws.Range("$col`5").Value = $row7t[$col]
    $ws.Range("$col`6").Value = $row5t[$col]
    $ws.Range("$col`7").Value = $row6t[$col]
}

$ws.Range("I5").Value = "'" + $i7
$ws.Range("I6").Value = "'" + $i5
$ws.Range("I7").Value = "'" + $i6

# row5 had no "Publik kommentar" before, so after the rotation row6 (which
# now holds the old row5 data) must be cleared once its old content has
# been moved to row5.
$ws.Range("AC5").Value = $ac7
if ($ac5 -eq $null -or $ac5 -eq "") {
    $ws.Range("AC6").ClearContents()
} else {
    $ws.Range("AC6").Value = $ac5
}
$ws.Range("AC7").Value = $ac6
